# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table (graphic frame / shape 2) gets its table style
#    switched from the custom "{E3DE31D3-...}" style to the built-in
#    style "{37FA73BC-80A3-41F2-8C43-6E8915A7F456}".
#
# 2) The deck's theme colour scheme ("Integral" / "Red Violet", currently
#    applied to the Slide Master / presentation theme) is switched back
#    to the stock "Office" colour palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{37FA73BC-80A3-41F2-8C43-6E8915A7F456}")

# --- 2. Theme colours -------------------------------------------------
function Convert-HexToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Index order for ThemeColorScheme.Item(n):
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToBgr $officeColors[$i - 1]
}
